# Refresh cryptocurrency Price (D) / Volume(1h) (E) columns with the
# latest coinranking.com scrape. Both columns are plain text cells in
# the workbook, so numeric-looking prices get a leading apostrophe to
# keep Excel from converting them to floating point numbers (which would
# lose trailing zeros / exact formatting, e.g. "191.70" -> 191.7).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '25.742.35'
$ws.Range("E2").Value = '  -0.24%  '

$ws.Range("D3").Value = '1.629.73'
$ws.Range("E3").Value = '  -0.48%  '

$ws.Range("D5").Value = "'214.36"
$ws.Range("E5").Value = '  -0.49%  '

$ws.Range("E6").Value = '  -0.80%  '

$ws.Range("E7").Value = '  -0.04%  '

$ws.Range("E8").Value = '  -0.91%  '

$ws.Range("E9").Value = '  -1.64%  '

$ws.Range("D10").Value = "'19.45"
$ws.Range("E10").Value = '  -2.30%  '

$ws.Range("D11").Value = "'0.0793"
$ws.Range("E11").Value = '  +0.79%  '

$ws.Range("E12").Value = '  +0.11%  '

$ws.Range("D13").Value = '1.855.58'
$ws.Range("E13").Value = '  -0.41%  '

$ws.Range("D14").Value = '1.629.55'
$ws.Range("E14").Value = '  -0.47%  '

$ws.Range("E15").Value = '  -0.01%  '

$ws.Range("E16").Value = '  -2.37%  '

$ws.Range("E17").Value = '  -0.18%  '

$ws.Range("D18").Value = '25.760.39'
$ws.Range("E18").Value = '  -0.28%  '

$ws.Range("E19").Value = '  -0.03%  '

$ws.Range("E20").Value = '  -0.28%  '

$ws.Range("D21").Value = "'191.70"
$ws.Range("E21").Value = '  -1.40%  '

$ws.Range("E22").Value = '  -0.37%  '

$ws.Range("E23").Value = '  +1.70%  '

$ws.Range("E24").Value = '  -0.04%  '

$ws.Range("D25").Value = "'1.82"
$ws.Range("E25").Value = '  +3.27%  '

$ws.Range("D26").Value = "'142.99"
$ws.Range("E26").Value = '  +1.92%  '

$ws.Range("E28").Value = '  +0.15%  '

$ws.Range("E29").Value = '  -0.76%  '

$ws.Range("E30").Value = '  -0.45%  '

$ws.Range("E31").Value = '  -0.91%  '

$ws.Range("E32").Value = '  +0.11%  '

$ws.Range("E33").Value = '  -1.17%  '

$ws.Range("D34").Value = "'1.57"
$ws.Range("E34").Value = '  -1.59%  '

$ws.Range("E35").Value = '  -0.27%  '

$ws.Range("E36").Value = '  +0.35%  '

$ws.Range("D37").Value = '1.133.01'
$ws.Range("E37").Value = '  +1.65%  '

$ws.Range("E38").Value = '  -1.99%  '

$ws.Range("D39").Value = "'0.542"
$ws.Range("E39").Value = '  -1.73%  '

$ws.Range("D40").Value = "'0.0155"
$ws.Range("E40").Value = '  -1.44%  '

$ws.Range("E41").Value = '  +0.13%  '

$ws.Range("E42").Value = '  +0.73%  '

$ws.Range("D43").Value = "'100.20"
$ws.Range("E43").Value = '  +1.04%  '

$ws.Range("E44").Value = '  -1.08%  '

$ws.Range("E45").Value = '  -0.23%  '

$ws.Range("D46").Value = '1.764.91'
$ws.Range("E46").Value = '  -0.28%  '

$ws.Range("E47").Value = '  +0.41%  '

$ws.Range("D48").Value = "'55.30"
$ws.Range("E48").Value = '  -0.53%  '

$ws.Range("E49").Value = '  +0.74%  '

$ws.Range("E50").Value = '  +0.11%  '

$ws.Range("D51").Value = "'1.41"
$ws.Range("E51").Value = '  +2.46%  '
